$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = [double]"0.02321097624268945"
$ws.Range("E3").Value = [double]"4.501101048536752e-10"
$ws.Range("E4").Value = [double]"7.98534541565203e-05"
$ws.Range("E5").Value = [double]"1.262106585610939e-09"
$ws.Range("E6").Value = [double]"2.48011777643617e-10"
$ws.Range("E7").Value = [double]"1.758038848441093e-05"
$ws.Range("E8").Value = [double]"1.46154786129859e-10"
$ws.Range("E9").Value = [double]"1.777843225227678e-05"
$ws.Range("E10").Value = [double]"0.01113417299502995"
$ws.Range("E11").Value = [double]"1.110178587469519e-07"
$ws.Range("E12").Value = [double]"5.414020516292887e-06"
$ws.Range("E13").Value = [double]"5.825431234064432e-05"
$ws.Range("E14").Value = [double]"0.0007935507406913846"
$ws.Range("E15").Value = [double]"0.03180380440829048"
$ws.Range("E16").Value = [double]"0.01243085060845515"
$ws.Range("E17").Value = [double]"6.789244634457656e-08"
